$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) from row 2 to row 23: 45184 -> 45186
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}
